$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at the top of the data (row 2), shifting all
#        existing price rows down by one. This also grows the used range
#        from A1:F185 to A1:F186 (the final row is just the old last row
#        shifted down, nothing needs to be appended at the bottom).
$ws.Rows.Item(2).Insert()

# --- 2. Populate the freshly inserted row 2 with the newest price entry.
#        (Same description/code/price/circular as the previous top row,
#        only the "Date" column advances by one day.)
$ws.Range("A2").Value = "07-02-2026"
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 320.45
$ws.Range("D2").NumberFormat = "0.000"
$ws.Range("E2").Value = "01-02-2026"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-02-2026.pdf"

# --- 3. The sheet's Hyperlinks collection does not automatically follow
#        rows that get shifted by Insert(), so every hyperlink ends up
#        anchored to the wrong cell. Rebuild the whole column from
#        scratch: every cell in F2:F186 links to the URL that is its own
#        displayed text.
$ws.Hyperlinks.Delete()

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("F$r")
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
}

# Hyperlinks.Add() re-styles the target cells with the built-in
# "Hyperlink" style (underline + theme color). Restore the original
# plain formatting by re-pasting the formats from column A (which keeps
# the untouched base style) back onto column F.
$ws.Range("A2:A$lastRow").Copy() | Out-Null
$ws.Range("F2:F$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "Done. Last row: $lastRow"
